$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace generic "Community N" labels with the real community names
# as found in the HW pdf document (first table: A5:A10, second table: A15:A20)
$names = @("DB (0)", "AL & ML (1)", "IR (2)", "DM (3)", "AL & TH (4)", "CV (5)")

for ($i = 0; $i -lt $names.Length; $i++) {
    $ws.Range("A" + (5 + $i)).Value = $names[$i]
    $ws.Range("A" + (15 + $i)).Value = $names[$i]
}

$ws.Range("N9").Select() | Out-Null
